$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.556007222541788
$ws.Range("C2").Value = 0.612606004275779
$ws.Range("K2").Value = 0.438282499020205
$ws.Range("L2").Value = 0.595756191953926
$ws.Range("N2").Value = 0.508231644030169

# Row 3
$ws.Range("B3").Value = 0.50312291438834
$ws.Range("K3").Value = 0.351050882580874
$ws.Range("L3").Value = 0.604399737467109
$ws.Range("N3").Value = 0.447300643788012

# Row 4
$ws.Range("B4").Value = 0.609601586795904
$ws.Range("K4").Value = 0.571096670838126
$ws.Range("L4").Value = 0.778963825426238
$ws.Range("N4").Value = 0.472900191628792

# Row 5
$ws.Range("B5").Value = 0.347853243460036
$ws.Range("N5").Value = 0.306445646731996

# Row 6
$ws.Range("B6").Value = 0.410626908494325
$ws.Range("K6").Value = 0.319383802321488
$ws.Range("L6").Value = 0.389233362357354
$ws.Range("N6").Value = 0.40055514051731
